# Generate Report for Handoff
#
# This mirrors a refreshed run of the localization-status report:
#  - The "Status" text moves from "Handed back: in sync with en-US" to
#    "Ready for handoff" on the Overview sheet (both language columns) and
#    on each language-specific sheet.
#  - The "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" stamps
#    are refreshed to the new generation time.
#  - The (now shorter) status column auto-shrinks to fit its new contents.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# --- Status text refresh -------------------------------------------------
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws2.Range("C2").Value = $newStatus
$ws3.Range("C2").Value = $newStatus

# --- Timestamp refresh -----------------------------------------------------
$ws1.Range("G2").Value = "2016-08-23 21:03:11"
$ws3.Range("H2").Value = "2016-08-23 21:03:11"
$ws2.Range("H2").Value = "2016-08-23 21:03:00"

# --- Column width shrink to fit the shorter status text ---------------------
# (ColumnWidth is expressed in characters and gets snapped to the engine's
# pixel grid on save; 16.3 is the input that lands closest to the target
# ~17.216-character rendered width.)
$ws1.Columns.Item(5).ColumnWidth = 16.3
$ws1.Columns.Item(6).ColumnWidth = 16.3
$ws2.Columns.Item(3).ColumnWidth = 16.3
$ws3.Columns.Item(3).ColumnWidth = 16.3
